$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 27
$ws.Range("A27").Value = 112529898
$ws.Range("B27").Value = 89571
$ws.Range("D27").Value = 'NT'
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = 'Granticka'
$ws.Range("G27").Value = 'Porodaedalea chrysoloma'
$ws.Range("H27").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q27").Value = 707896
$ws.Range("R27").Value = 7087596
$ws.Range("AX27").Value = 'Simon Mattsson, Maja Östlund'

# Row 28
$ws.Range("A28").Value = 112529839
$ws.Range("B28").Value = 89553
$ws.Range("D28").Value = 'NT'
$ws.Range("E28").Value = 1202
$ws.Range("F28").Value = 'Ullticka'
$ws.Range("G28").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H28").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q28").Value = 708114
$ws.Range("R28").Value = 7087436
$ws.Range("AX28").Value = 'Simon Mattsson, Åsa Stenman'

# Row 29
$ws.Range("A29").Value = 112529887
$ws.Range("B29").Value = 89517
$ws.Range("D29").Value = 'LC'
$ws.Range("E29").Value = 5447
$ws.Range("F29").Value = 'Vedticka'
$ws.Range("G29").Value = 'Fuscoporia viticola'
$ws.Range("H29").Value = '(Schwein.) Murrill'
$ws.Range("Q29").Value = 707931
$ws.Range("R29").Value = 7087583
$ws.Range("AX29").Value = 'Simon Mattsson'

# Row 30
$ws.Range("A30").Value = 112529871
$ws.Range("B30").Value = 89571
$ws.Range("D30").Value = 'NT'
$ws.Range("E30").Value = 5432
$ws.Range("F30").Value = 'Granticka'
$ws.Range("G30").Value = 'Porodaedalea chrysoloma'
$ws.Range("H30").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q30").Value = 708172
$ws.Range("R30").Value = 7087543
$ws.Range("AX30").Value = 'Simon Mattsson, Maja Östlund'

# Row 31
$ws.Range("A31").Value = 112529879
$ws.Range("B31").Value = 89553
$ws.Range("D31").Value = 'NT'
$ws.Range("E31").Value = 1202
$ws.Range("F31").Value = 'Ullticka'
$ws.Range("G31").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H31").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q31").Value = 707931
$ws.Range("R31").Value = 7087576
$ws.Range("AX31").Value = 'Simon Mattsson, Åsa Stenman'

# Row 32
$ws.Range("A32").Value = 112529842
$ws.Range("B32").Value = 96735
$ws.Range("D32").Value = 'VU'
$ws.Range("E32").Value = 220787
$ws.Range("F32").Value = 'Knärot'
$ws.Range("G32").Value = 'Goodyera repens'
$ws.Range("H32").Value = '(L.) R. Br.'
$ws.Range("Q32").Value = 708088
$ws.Range("R32").Value = 7087457
$ws.Range("AX32").Value = 'Simon Mattsson, Maja Östlund'

# Row 33
$ws.Range("A33").Value = 112529844
$ws.Range("B33").Value = 96735
$ws.Range("D33").Value = 'VU'
$ws.Range("E33").Value = 220787
$ws.Range("F33").Value = 'Knärot'
$ws.Range("G33").Value = 'Goodyera repens'
$ws.Range("H33").Value = '(L.) R. Br.'
$ws.Range("Q33").Value = 708081
$ws.Range("R33").Value = 7087468
$ws.Range("AX33").Value = 'Simon Mattsson, Anna Hallmén'

# Row 34
$ws.Range("A34").Value = 112529872
$ws.Range("B34").Value = 96735
$ws.Range("D34").Value = 'VU'
$ws.Range("E34").Value = 220787
$ws.Range("F34").Value = 'Knärot'
$ws.Range("G34").Value = 'Goodyera repens'
$ws.Range("H34").Value = '(L.) R. Br.'
$ws.Range("Q34").Value = 708187
$ws.Range("R34").Value = 7087543
$ws.Range("AX34").Value = 'Simon Mattsson, Maja Östlund'

# Row 48
$ws.Range("A48").Value = 112529890
$ws.Range("B48").Value = 89993
$ws.Range("D48").Value = 'VU'
$ws.Range("E48").Value = 1209
$ws.Range("F48").Value = 'Rynkskinn'
$ws.Range("G48").Value = 'Phlebia centrifuga'
$ws.Range("H48").Value = 'P.Karst.'
$ws.Range("Q48").Value = 707943
$ws.Range("R48").Value = 7087588
$ws.Range("AX48").Value = 'Simon Mattsson, Åsa Stenman'

# Row 49
$ws.Range("A49").Value = 112529858
$ws.Range("B49").Value = 78647
$ws.Range("D49").Value = 'LC'
$ws.Range("E49").Value = 6456
$ws.Range("F49").Value = 'Skinnlav'
$ws.Range("G49").Value = 'Leptogium saturninum'
$ws.Range("H49").Value = '(Dicks.) Nyl.'
$ws.Range("Q49").Value = 708059
$ws.Range("R49").Value = 7087503
$ws.Range("AX49").Value = 'Simon Mattsson, Maja Östlund'

# Row 50
$ws.Range("A50").Value = 112529860
$ws.Range("B50").Value = 78746
$ws.Range("D50").Value = 'LC'
$ws.Range("E50").Value = 6463
$ws.Range("F50").Value = 'Bårdlav'
$ws.Range("G50").Value = 'Nephroma parile'
$ws.Range("H50").Value = '(Ach.) Ach.'
$ws.Range("Q50").Value = 708056
$ws.Range("R50").Value = 7087508
$ws.Range("AX50").Value = 'Simon Mattsson'

# Row 51
$ws.Range("A51").Value = 112529906
$ws.Range("B51").Value = 96735
$ws.Range("D51").Value = 'VU'
$ws.Range("E51").Value = 220787
$ws.Range("F51").Value = 'Knärot'
$ws.Range("G51").Value = 'Goodyera repens'
$ws.Range("H51").Value = '(L.) R. Br.'
$ws.Range("Q51").Value = 708230
$ws.Range("R51").Value = 7087674
$ws.Range("AX51").Value = 'Simon Mattsson'

# Row 52
$ws.Range("A52").Value = 112529907
$ws.Range("B52").Value = 96735
$ws.Range("D52").Value = 'VU'
$ws.Range("E52").Value = 220787
$ws.Range("F52").Value = 'Knärot'
$ws.Range("G52").Value = 'Goodyera repens'
$ws.Range("H52").Value = '(L.) R. Br.'
$ws.Range("Q52").Value = 708225
$ws.Range("R52").Value = 7087689
$ws.Range("AX52").Value = 'Simon Mattsson, Maja Östlund'

# Row 54
$ws.Range("A54").Value = 112529903
$ws.Range("B54").Value = 89564
$ws.Range("D54").Value = 'LC'
$ws.Range("E54").Value = 1205
$ws.Range("F54").Value = 'Stor aspticka'
$ws.Range("G54").Value = 'Phellinus populicola'
$ws.Range("H54").Value = 'Niemelä'
$ws.Range("Q54").Value = 708216
$ws.Range("R54").Value = 7087659
$ws.Range("AX54").Value = 'Simon Mattsson'

# Row 55
$ws.Range("A55").Value = 112529875
$ws.Range("B55").Value = 89553
$ws.Range("D55").Value = 'NT'
$ws.Range("E55").Value = 1202
$ws.Range("F55").Value = 'Ullticka'
$ws.Range("G55").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H55").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q55").Value = 707988
$ws.Range("R55").Value = 7087564
$ws.Range("AX55").Value = 'Simon Mattsson, Åsa Stenman'

# Row 56
$ws.Range("A56").Value = 112529854
$ws.Range("B56").Value = 96735
$ws.Range("D56").Value = 'VU'
$ws.Range("E56").Value = 220787
$ws.Range("F56").Value = 'Knärot'
$ws.Range("G56").Value = 'Goodyera repens'
$ws.Range("H56").Value = '(L.) R. Br.'
$ws.Range("Q56").Value = 708149
$ws.Range("R56").Value = 7087493
$ws.Range("AX56").Value = 'Simon Mattsson, Maja Östlund'

# Row 57
$ws.Range("A57").Value = 112529853
$ws.Range("B57").Value = 96735
$ws.Range("D57").Value = 'VU'
$ws.Range("E57").Value = 220787
$ws.Range("F57").Value = 'Knärot'
$ws.Range("G57").Value = 'Goodyera repens'
$ws.Range("H57").Value = '(L.) R. Br.'
$ws.Range("Q57").Value = 708163
$ws.Range("R57").Value = 7087493
$ws.Range("AX57").Value = 'Simon Mattsson'

# Row 58
$ws.Range("A58").Value = 112529909
$ws.Range("B58").Value = 90814
$ws.Range("D58").Value = 'LC'
$ws.Range("E58").Value = 4364
$ws.Range("F58").Value = 'Dropptaggsvamp'
$ws.Range("G58").Value = 'Hydnellum ferrugineum'
$ws.Range("H58").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q58").Value = 708221
$ws.Range("R58").Value = 7087718
$ws.Range("AX58").Value = 'Simon Mattsson, Åsa Stenman'

# Row 59
$ws.Range("A59").Value = 112529885
$ws.Range("B59").Value = 89571
$ws.Range("D59").Value = 'NT'
$ws.Range("E59").Value = 5432
$ws.Range("F59").Value = 'Granticka'
$ws.Range("G59").Value = 'Porodaedalea chrysoloma'
$ws.Range("H59").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q59").Value = 707987
$ws.Range("R59").Value = 7087580
$ws.Range("AX59").Value = 'Simon Mattsson, Maja Östlund'

# Row 75
$ws.Range("A75").Value = 112529901
$ws.Range("B75").Value = 89517
$ws.Range("D75").Value = 'LC'
$ws.Range("E75").Value = 5447
$ws.Range("F75").Value = 'Vedticka'
$ws.Range("G75").Value = 'Fuscoporia viticola'
$ws.Range("H75").Value = '(Schwein.) Murrill'
$ws.Range("Q75").Value = 708195
$ws.Range("R75").Value = 7087612
$ws.Range("AX75").Value = 'Simon Mattsson, Åsa Stenman'

# Row 76
$ws.Range("A76").Value = 112529876
$ws.Range("B76").Value = 89553
$ws.Range("D76").Value = 'NT'
$ws.Range("E76").Value = 1202
$ws.Range("F76").Value = 'Ullticka'
$ws.Range("G76").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H76").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q76").Value = 707983
$ws.Range("R76").Value = 7087566
$ws.Range("AX76").Value = 'Simon Mattsson, Maja Östlund'

# Row 77
$ws.Range("A77").Value = 112529893
$ws.Range("B77").Value = 89993
$ws.Range("D77").Value = 'VU'
$ws.Range("E77").Value = 1209
$ws.Range("F77").Value = 'Rynkskinn'
$ws.Range("G77").Value = 'Phlebia centrifuga'
$ws.Range("H77").Value = 'P.Karst.'
$ws.Range("Q77").Value = 707938
$ws.Range("R77").Value = 7087590
$ws.Range("AX77").Value = 'Simon Mattsson, Åsa Stenman'

# Row 78
$ws.Range("A78").Value = 112529862
$ws.Range("B78").Value = 78713
$ws.Range("D78").Value = 'NT'
$ws.Range("E78").Value = 6458
$ws.Range("F78").Value = 'Lunglav'
$ws.Range("G78").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H78").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q78").Value = 708051
$ws.Range("R78").Value = 7087517
$ws.Range("AX78").Value = 'Simon Mattsson, Åsa Stenman'
